$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full credit-rating scale (best -> worst), used to translate the textual
# "rating" column (E) into a numeric "rank" column (F).
$scale = @("AAA","AA+","AA","AA-","A+","A","A-","BBB+","BBB","BBB-","BB+","BB","BB-","B+","B","B-","C+","C","C-","D")

$rankOf = @{}
for ($i = 0; $i -lt $scale.Length; $i++) {
    $rankOf[$scale[$i]] = $i + 1
}

# New header for column F.
$ws.Range("F1").Value = "rank"

# Find last used row (column A) and loop row by row, looking up the rating
# text in column E and writing the corresponding numeric rank into F.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($i = 2; $i -le $lastRow; $i++) {
    $rating = $ws.Cells.Item($i, 5).Value2
    $rank = $rankOf[$rating]
    $ws.Cells.Item($i, 6).Value = $rank
}

# Mirror the saved view state from the edit session: scrolled down with the
# newly-added F37 cell selected.
$ws.Range("F37").Select()
$excel.ActiveWindow.ScrollRow = 22
